$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "iaest-measure:residencia-comarca-nombre"
$ws.Range("I2").Value = "iaest-measure:residencia-provincia-nombre"
$ws.Range("J2").Value = "iaest-measure:nacionalidad-area-nombre"
$ws.Range("K2").Value = "iaest-measure:edad-grandes-grupos"

$ws.Range("G4").Value = "xsd:int"
$ws.Range("I4").Value = "xsd:int"
$ws.Range("J4").Value = "xsd:int"
$ws.Range("K4").Value = "xsd:int"

$ws.Range("J5:K5").Clear()
